$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: update existing timestamp ---
$ws.Range("A2").Value = 43298.498043981483

# --- New shared-string text used by rows 11 & 12 ---
$abhishek = "Attendance Processed for Abhishek"

# --- Row 6 ---
$ws.Range("A2").Copy($ws.Range("A6"))
$ws.Range("A6").Value = 43298.469664351855
$ws.Range("C6").Value = "Attendance Processed for BatchA"

# --- Row 7 ---
$ws.Range("A2").Copy($ws.Range("A7"))
$ws.Range("A7").Value = 43298.46980324074
$ws.Range("C7").Value = "Attendance Processed for BatchB"

# --- Row 8 (wrapped long error message, tall row) ---
$ws.Range("A2").Copy($ws.Range("A8"))
$ws.Range("A8").Value = 43298.469814814816
$ws.Range("C5").Copy($ws.Range("C8"))
$ws.Rows.Item(8).RowHeight = 409.6

# --- Row 9 ---
$ws.Range("A2").Copy($ws.Range("A9"))
$ws.Range("A9").Value = 43298.477280092593
$ws.Range("C9").Value = "Attendance Processed for BatchA"

# --- Row 10 ---
$ws.Range("A2").Copy($ws.Range("A10"))
$ws.Range("A10").Value = 43298.477465277778
$ws.Range("C10").Value = "Attendance Processed for BatchB"

# --- Row 11 ---
$ws.Range("A2").Copy($ws.Range("A11"))
$ws.Range("A11").Value = 43298.479513888888
$ws.Range("C11").Value = $abhishek

# --- Row 12 ---
$ws.Range("A2").Copy($ws.Range("A12"))
$ws.Range("A12").Value = 43298.481631944444
$ws.Range("C12").Value = $abhishek

# --- Row 13 ---
$ws.Range("A2").Copy($ws.Range("A13"))
$ws.Range("A13").Value = 43298.481840277775
$ws.Range("C13").Value = "Attendance Processed for BatchB"

# --- Row 14 ---
$ws.Range("A2").Copy($ws.Range("A14"))
$ws.Range("A14").Value = 43298.49832175926
$ws.Range("C14").Value = "Attendance Processed for BatchA"

# --- Row 15 ---
$ws.Range("A2").Copy($ws.Range("A15"))
$ws.Range("A15").Value = 43298.498379629629
$ws.Range("C15").Value = "Attendance Processed for BatchB"
